$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$status = "Handed back: in sync with en-US"

# Overview sheet: row 3 corresponds to c2234ac1 file
$overview.Range("B3").Value = $status
$overview.Range("C3").Value = $status

# zh-cn sheet: row 3 corresponds to c2234ac1 file
$zhcn.Range("C3").Value = $status
$zhcn.Range("H3").Value = "2016-03-14 09:18:36"

# de-de sheet: row 3 corresponds to c2234ac1 file
$dede.Range("C3").Value = $status
$dede.Range("H3").Value = "2016-03-14 09:18:49"
